$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 197, shifting existing rows 197:219 down to 198:220
$ws.Rows.Item(197).Insert()

# Populate the new row 197 with the new price observation
$ws.Cells.Item(197, 1).Value = 11
$ws.Cells.Item(197, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(197, 3).Value = "Bíobío"
$ws.Cells.Item(197, 4).Value = 44504
$ws.Cells.Item(197, 5).Value = 8
$ws.Cells.Item(197, 6).Value = 100114014
$ws.Cells.Item(197, 7).Value = "Betarraga"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 1700
$ws.Cells.Item(197, 11).Value = 600
$ws.Cells.Item(197, 12).Value = 650
$ws.Cells.Item(197, 13).Value = 624
$ws.Cells.Item(197, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(197, 15).Value = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value = 125
$ws.Cells.Item(197, 17).Value = 5
$ws.Cells.Item(197, 18).Value = "Hortaliza"
